$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.854.81"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.39"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.33"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4827"
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07370"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9405"
$ws.Range("E10").Value = "  +0.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.05"
$ws.Range("E11").Value = "  +4.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07790"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.885.02"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.519"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.617"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.37"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008850"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.876.06"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.86"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.123"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.122.45"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.947"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.51"
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.58"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.92"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.974"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08888"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.344"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.224"
$ws.Range("E33").Value = "  +3.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7707"
$ws.Range("E34").Value = "  +4.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.657"
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.729"
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.128"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02046"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5628"
$ws.Range("E39").Value = "  +5.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05373"
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.000"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.065"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.552"
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1529"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4876"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.64"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.34"
$ws.Range("E47").Value = "  +2.84%  "
$ws.Range("E48").Value = "  -0.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.665"
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.09"
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06123"
$ws.Range("E51").Value = "  +0.78%  "
